$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.716.35'
$ws.Range('E2').Value = '  -4.53%  '
$ws.Range('D3').Value = '2.449.83'
$ws.Range('E3').Value = '  -6.00%  '
$ws.Range('E4').Value = '  +0.00%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '546.68'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -4.29%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '144.29'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -6.90%  '
$ws.Range('E7').Value = '  -0.01%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.594'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -4.25%  '
$ws.Range('D9').Value = '2.447.88'
$ws.Range('E9').Value = '  -5.97%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.106'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -8.53%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.153'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -2.01%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '5.35'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -7.99%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.350'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -7.53%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '25.86'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -7.22%  '
$ws.Range('D15').Value = '2.893.41'
$ws.Range('E15').Value = '  -5.90%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.0000162'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -9.51%  '
$ws.Range('D17').Value = '60.655.72'
$ws.Range('E17').Value = '  -4.46%  '
$ws.Range('D18').Value = '2.457.75'
$ws.Range('E18').Value = '  -5.64%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '10.99'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -7.76%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '6.88'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -7.98%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '4.15'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -7.45%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '317.29'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -6.82%  '
$ws.Range('E23').Value = '  -0.19%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '63.23'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -5.81%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '1.76'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -2.61%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.578.92'
$ws.Range('E26').Value = '  -5.66%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0₃0966'
$ws.Range('E27').Value = '  -9.07%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '532.75'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -7.89%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '1.47'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -4.96%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '8.29'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -8.41%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '7.51'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -3.61%  '
$ws.Range('E33').Value = '  -7.25%  '
$ws.Range('E34').Value = '  -7.93%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.56'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -9.15%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '5.81'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -11.62%  '
$ws.Range('E37').Value = '  -0.05%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '4.79'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -9.93%  '
$ws.Range('E39').Value = '  -6.36%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '18.38'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -6.17%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '145.64'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -5.35%  '
$ws.Range('E42').Value = '  +0.02%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '1.69'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -9.20%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '39.77'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -4.27%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -9.32%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '145.90'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -6.94%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '3.54'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -8.18%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '20.68'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -11.91%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.0526'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -9.50%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.580'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -7.47%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.0935'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -6.26%  '
